# Async ajax for save filter
# Adds two new timeline entries (rows 23 and 24) that were previously
# blank placeholder rows, and moves the active selection to A23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (C22) already carries the date number-format style used
# throughout column C; copy that formatting onto C23:C24 before writing
# the new date values so they pick up the same style index instead of
# minting a new custom number format.
$ws.Range("C22").Copy()
$ws.Range("C23:C24").PasteSpecial(-4122)

# New shared strings must be appended in the same order the original
# workbook used (index 22 = "Запись фильтров", index 23 = "Очистка
# фильтров..."), so write row 24's text before row 23's.

# Row 24: "Запись фильтров" -> becomes shared string index 22
$ws.Range("A24").Value = "Работа по созданию функционала загрузки данных (Запись фильтров)"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 43570

# Row 23: "Очистка фильтров, Таблица фильтров, модель хранения" -> index 23
$ws.Range("A23").Value = "Работа по созданию функционала загрузки данных (Очистка фильтров, Таблица фильтров, модель хранения)"
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 43569

# Column D keeps its existing shared formula (B*$B$1) and recalculates
# automatically; the grand-total formula in D38 follows suit.

# Move/restore the active selection to A23, matching the saved view state.
$ws.Range("A23").Select()
